# Append a new data row (row 98) to "Tabela1" on sheet1, mirroring the
# formatting of the last existing data row (row 97), fill in the new
# day's COVID-19 stats, and grow the table / autofilter to include it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 97
$newRow = 98

# Duplicate row 97's formatting into the new row 98 (keeps the same
# cellXf/style indices the sheet already uses for data rows).
$ws.Rows($lastRow).Copy()
$ws.Rows($newRow).Insert()

# New day's values: Date, Tested(all), Tested(daily), Positive(all),
# Positive(daily), Hospitalized, IntensiveCare, Discharged, Deaths(all),
# Deaths(daily)
$values = @(43998, 89151, 986, 1503, 4, 7, 1, 0, 109, 0)

for ($col = 1; $col -le 10; $col++) {
    $ws.Cells.Item($newRow, $col).Value = $values[$col - 1]
}

# Grow the table (and its autofilter range) to cover the new row.
$table = $ws.ListObjects.Item("Tabela1")
$newRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($newRow, 10))
$table.Resize($newRange)

# Mirror the post-edit selection state (new last row selected).
$ws.Range("A98:J98").Select()
